$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Weekly roll: rows 191-219 shift D/J/K/L/M/P down by one week ---
# Row 191
$ws.Cells.Item(191, 4).Value = 44505
$ws.Cells.Item(191, 10).Value = 760
$ws.Cells.Item(191, 11).Value = 6500
$ws.Cells.Item(191, 12).Value = 7000
$ws.Cells.Item(191, 13).Value = 6750
$ws.Cells.Item(191, 16).Value = 338

# Row 192
$ws.Cells.Item(192, 4).Value = 44487
$ws.Cells.Item(192, 10).Value = 600
$ws.Cells.Item(192, 11).Value = 6500
$ws.Cells.Item(192, 12).Value = 7000
$ws.Cells.Item(192, 13).Value = 6750
$ws.Cells.Item(192, 16).Value = 338

# Row 193
$ws.Cells.Item(193, 4).Value = 44425
$ws.Cells.Item(193, 10).Value = 660
$ws.Cells.Item(193, 11).Value = 5000
$ws.Cells.Item(193, 12).Value = 5500
$ws.Cells.Item(193, 13).Value = 5250
$ws.Cells.Item(193, 16).Value = 262

# Row 194
$ws.Cells.Item(194, 4).Value = 44343
$ws.Cells.Item(194, 10).Value = 700
$ws.Cells.Item(194, 11).Value = 4800
$ws.Cells.Item(194, 12).Value = 5000
$ws.Cells.Item(194, 13).Value = 4900
$ws.Cells.Item(194, 16).Value = 245

# Row 195
$ws.Cells.Item(195, 4).Value = 44370
$ws.Cells.Item(195, 10).Value = 800
$ws.Cells.Item(195, 11).Value = 5000
$ws.Cells.Item(195, 12).Value = 5500
$ws.Cells.Item(195, 13).Value = 5250
$ws.Cells.Item(195, 16).Value = 262

# Row 196
$ws.Cells.Item(196, 4).Value = 44449
$ws.Cells.Item(196, 10).Value = 700
$ws.Cells.Item(196, 11).Value = 4500
$ws.Cells.Item(196, 12).Value = 5000
$ws.Cells.Item(196, 13).Value = 4750
$ws.Cells.Item(196, 16).Value = 238

# Row 197
$ws.Cells.Item(197, 4).Value = 44168
$ws.Cells.Item(197, 10).Value = 600
$ws.Cells.Item(197, 11).Value = 5000
$ws.Cells.Item(197, 12).Value = 5500
$ws.Cells.Item(197, 13).Value = 5250
$ws.Cells.Item(197, 16).Value = 262

# Row 198
$ws.Cells.Item(198, 4).Value = 44175
$ws.Cells.Item(198, 10).Value = 600
$ws.Cells.Item(198, 11).Value = 5000
$ws.Cells.Item(198, 12).Value = 5500
$ws.Cells.Item(198, 13).Value = 5250
$ws.Cells.Item(198, 16).Value = 262

# Row 199
$ws.Cells.Item(199, 4).Value = 44392
$ws.Cells.Item(199, 10).Value = 700
$ws.Cells.Item(199, 11).Value = 5000
$ws.Cells.Item(199, 12).Value = 5500
$ws.Cells.Item(199, 13).Value = 5250
$ws.Cells.Item(199, 16).Value = 262

# Row 200
$ws.Cells.Item(200, 4).Value = 44286
$ws.Cells.Item(200, 10).Value = 800
$ws.Cells.Item(200, 11).Value = 5000
$ws.Cells.Item(200, 12).Value = 5500
$ws.Cells.Item(200, 13).Value = 5250
$ws.Cells.Item(200, 16).Value = 262

# Row 201
$ws.Cells.Item(201, 4).Value = 44473
$ws.Cells.Item(201, 10).Value = 600
$ws.Cells.Item(201, 11).Value = 6000
$ws.Cells.Item(201, 12).Value = 7000
$ws.Cells.Item(201, 13).Value = 6500
$ws.Cells.Item(201, 16).Value = 325

# Row 202
$ws.Cells.Item(202, 4).Value = 44400
$ws.Cells.Item(202, 10).Value = 720
$ws.Cells.Item(202, 11).Value = 5000
$ws.Cells.Item(202, 12).Value = 5500
$ws.Cells.Item(202, 13).Value = 5250
$ws.Cells.Item(202, 16).Value = 262

# Row 203
$ws.Cells.Item(203, 4).Value = 44484
$ws.Cells.Item(203, 10).Value = 760
$ws.Cells.Item(203, 11).Value = 6500
$ws.Cells.Item(203, 12).Value = 7000
$ws.Cells.Item(203, 13).Value = 6750
$ws.Cells.Item(203, 16).Value = 338

# Row 204
$ws.Cells.Item(204, 4).Value = 44181
$ws.Cells.Item(204, 10).Value = 400
$ws.Cells.Item(204, 11).Value = 5000
$ws.Cells.Item(204, 12).Value = 5500
$ws.Cells.Item(204, 13).Value = 5250
$ws.Cells.Item(204, 16).Value = 262

# Row 205
$ws.Cells.Item(205, 4).Value = 44494
$ws.Cells.Item(205, 10).Value = 600
$ws.Cells.Item(205, 11).Value = 6800
$ws.Cells.Item(205, 12).Value = 7000
$ws.Cells.Item(205, 13).Value = 6900
$ws.Cells.Item(205, 16).Value = 345

# Row 206
$ws.Cells.Item(206, 4).Value = 44342
$ws.Cells.Item(206, 10).Value = 800
$ws.Cells.Item(206, 11).Value = 4800
$ws.Cells.Item(206, 12).Value = 5000
$ws.Cells.Item(206, 13).Value = 4900
$ws.Cells.Item(206, 16).Value = 245

# Row 207
$ws.Cells.Item(207, 4).Value = 44328
$ws.Cells.Item(207, 10).Value = 800
$ws.Cells.Item(207, 11).Value = 4800
$ws.Cells.Item(207, 12).Value = 5000
$ws.Cells.Item(207, 13).Value = 4900
$ws.Cells.Item(207, 16).Value = 245

# Row 208
$ws.Cells.Item(208, 4).Value = 44301
$ws.Cells.Item(208, 10).Value = 700
$ws.Cells.Item(208, 11).Value = 5000
$ws.Cells.Item(208, 12).Value = 5500
$ws.Cells.Item(208, 13).Value = 5250
$ws.Cells.Item(208, 16).Value = 262

# Row 209
$ws.Cells.Item(209, 4).Value = 44330
$ws.Cells.Item(209, 10).Value = 800
$ws.Cells.Item(209, 11).Value = 4800
$ws.Cells.Item(209, 12).Value = 5000
$ws.Cells.Item(209, 13).Value = 4900
$ws.Cells.Item(209, 16).Value = 245

# Row 210
$ws.Cells.Item(210, 4).Value = 44270
$ws.Cells.Item(210, 10).Value = 600
$ws.Cells.Item(210, 11).Value = 5500
$ws.Cells.Item(210, 12).Value = 6000
$ws.Cells.Item(210, 13).Value = 5750
$ws.Cells.Item(210, 16).Value = 288

# Row 211
$ws.Cells.Item(211, 4).Value = 44295
$ws.Cells.Item(211, 10).Value = 800
$ws.Cells.Item(211, 11).Value = 5000
$ws.Cells.Item(211, 12).Value = 5500
$ws.Cells.Item(211, 13).Value = 5250
$ws.Cells.Item(211, 16).Value = 262

# Row 212
$ws.Cells.Item(212, 4).Value = 44217
$ws.Cells.Item(212, 10).Value = 600
$ws.Cells.Item(212, 11).Value = 5000
$ws.Cells.Item(212, 12).Value = 5500
$ws.Cells.Item(212, 13).Value = 5250
$ws.Cells.Item(212, 16).Value = 262

# Row 213
$ws.Cells.Item(213, 4).Value = 44421
$ws.Cells.Item(213, 10).Value = 700
$ws.Cells.Item(213, 11).Value = 5000
$ws.Cells.Item(213, 12).Value = 5500
$ws.Cells.Item(213, 13).Value = 5250
$ws.Cells.Item(213, 16).Value = 262

# Row 214
$ws.Cells.Item(214, 4).Value = 44383
$ws.Cells.Item(214, 10).Value = 600
$ws.Cells.Item(214, 11).Value = 5000
$ws.Cells.Item(214, 12).Value = 5500
$ws.Cells.Item(214, 13).Value = 5250
$ws.Cells.Item(214, 16).Value = 262

# Row 215
$ws.Cells.Item(215, 4).Value = 44244
$ws.Cells.Item(215, 10).Value = 800
$ws.Cells.Item(215, 11).Value = 5500
$ws.Cells.Item(215, 12).Value = 6000
$ws.Cells.Item(215, 13).Value = 5750
$ws.Cells.Item(215, 16).Value = 288

# Row 216
$ws.Cells.Item(216, 4).Value = 44307
$ws.Cells.Item(216, 10).Value = 800
$ws.Cells.Item(216, 11).Value = 5000
$ws.Cells.Item(216, 12).Value = 5500
$ws.Cells.Item(216, 13).Value = 5250
$ws.Cells.Item(216, 16).Value = 262

# Row 217
$ws.Cells.Item(217, 4).Value = 44273
$ws.Cells.Item(217, 10).Value = 700
$ws.Cells.Item(217, 11).Value = 6000
$ws.Cells.Item(217, 12).Value = 6500
$ws.Cells.Item(217, 13).Value = 6250
$ws.Cells.Item(217, 16).Value = 312

# Row 218
$ws.Cells.Item(218, 4).Value = 44433
$ws.Cells.Item(218, 10).Value = 800
$ws.Cells.Item(218, 11).Value = 5000
$ws.Cells.Item(218, 12).Value = 5500
$ws.Cells.Item(218, 13).Value = 5250
$ws.Cells.Item(218, 16).Value = 262

# Row 219
$ws.Cells.Item(219, 4).Value = 44302
$ws.Cells.Item(219, 10).Value = 800
$ws.Cells.Item(219, 11).Value = 5000
$ws.Cells.Item(219, 12).Value = 5500
$ws.Cells.Item(219, 13).Value = 5250
$ws.Cells.Item(219, 16).Value = 262

# --- New row 220: duplicate of original (pre-edit) row 219 ---
$ws.Cells.Item(220, 1).Value = 8
$ws.Cells.Item(220, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(220, 3).Value = "Coquimbo"
$ws.Cells.Item(220, 4).Value = 44179
$ws.Cells.Item(220, 5).Value = 4
$ws.Cells.Item(220, 6).Value = 100114013
$ws.Cells.Item(220, 7).Value = "Zanahoria"
$ws.Cells.Item(220, 8).Value = "Sin especificar"
$ws.Cells.Item(220, 9).Value = "Primera"
$ws.Cells.Item(220, 10).Value = 760
$ws.Cells.Item(220, 11).Value = 5000
$ws.Cells.Item(220, 12).Value = 5500
$ws.Cells.Item(220, 13).Value = 5250
$ws.Cells.Item(220, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(220, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(220, 16).Value = 262
$ws.Cells.Item(220, 17).Value = 20
$ws.Cells.Item(220, 18).Value = "Hortaliza"
$ws.Cells.Item(220, 4).NumberFormat = $ws.Cells.Item(219, 4).NumberFormat
